$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83-128 down to 84-129.
$ws.Rows("83:83").Insert()

# Fill the new row 83 with the new record's data (constant columns shared
# with all other rows in this sheet, plus the new date/price values).
$ws.Range("A83").Value2 = 8
$ws.Range("B83").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C83").Value2 = "Coquimbo"
$ws.Range("D83").Value2 = 44510
$ws.Range("E83").Value2 = 4
$ws.Range("F83").Value2 = 100112037
$ws.Range("G83").Value2 = "Cebollín"
$ws.Range("H83").Value2 = "Sin especificar"
$ws.Range("I83").Value2 = "Primera"
$ws.Range("J83").Value2 = 3200
$ws.Range("K83").Value2 = 900
$ws.Range("L83").Value2 = 1000
$ws.Range("M83").Value2 = 950
$ws.Range("N83").Value2 = '$/paquete 6 unidades'
$ws.Range("O83").Value2 = "Provincia del Elquí"
$ws.Range("P83").Value2 = 158
$ws.Range("Q83").Value2 = 6
$ws.Range("R83").Value2 = "Hortaliza"
